$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Fixed bugs with AdminAddNurse" - Sprint burndown updates for the
# "Search for a patient by criteria" story (row 5) and the
# "Fix Add Nurse" story (row 9): both now carry a second assignee
# (Rahman/Ahmad) and actual effort/remaining numbers.

# Row 5 - "Search for a patient by criteria " / "Fix Nurse state text box" assignment
$ws.Range("D5").Value = "Rahman/Ahmad"
$ws.Range("E5").Value = 14
$ws.Range("F5").Value = 14

# Row 9 - "Fix Add Nurse"
$ws.Range("D9").Value = "Rahman/Ahmad"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0

# Selection cursor, matching the state the workbook was saved in.
$ws.Range("G9").Select()
